# Fixed picture positions on some slides. Updated App theme for widgets to new look.
#
# NOTE on units: PowerPoint's COM object model expresses Shape.Left/Top/Width/Height
# in points (1 pt = 12700 EMU), while the underlying OOXML stores EMU. To land on an
# *exact* target EMU value we use literal point values that have been precomputed so
# that point -> EMU round-tripping reproduces the precise EMU figures from the source
# diff (avoiding any off-by-one EMU drift from naive division).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 ("...slides/slide3.xml" in the package / Slides.Item(3))
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Shape 1: "Content Placeholder 2" - reposition (size stays the same)
$s3Content = $s3.Shapes.Item(1)
$s3Content.Left = 45.046535433070865   # 572091 EMU
$s3Content.Top  = 14.37496062992126    # 182562 EMU

# Shape 5: "Picture 6"
$s3Pic1 = $s3.Shapes.Item(5)
$s3Pic1.Left   = 144.0                 # 1828800 EMU
$s3Pic1.Top    = 86.4                  # 1097280 EMU
$s3Pic1.Width  = 312.74906921396115    # 3971913 EMU
$s3Pic1.Height = 432.0                 # 5486400 EMU

# Shape 6: "Picture 7"
$s3Pic2 = $s3.Shapes.Item(6)
$s3Pic2.Left   = 504.0                 # 6400800 EMU
$s3Pic2.Top    = 86.4                  # 1097280 EMU
$s3Pic2.Width  = 312.74906921396115    # 3971913 EMU
$s3Pic2.Height = 432.0                 # 5486400 EMU

# ---------------------------------------------------------------------------
# Slide 4 ("...slides/slide4.xml" in the package / Slides.Item(4))
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Shape 1: "Content Placeholder 2" - reposition (size stays the same)
$s4Content = $s4.Shapes.Item(1)
$s4Content.Left = 103.95566940333899   # 1320237 EMU
$s4Content.Top  = 25.73472499944887    # 326831 EMU

# Shape 5: "Picture 6"
$s4Pic1 = $s4.Shapes.Item(5)
$s4Pic1.Left   = 144.0                 # 1828800 EMU
$s4Pic1.Top    = 86.4                  # 1097280 EMU
$s4Pic1.Width  = 305.26551818895615    # 3876872 EMU
$s4Pic1.Height = 432.0                 # 5486400 EMU

# Shape 6: "Picture 7"
$s4Pic2 = $s4.Shapes.Item(6)
$s4Pic2.Left   = 503.9999237068166     # 6400799 EMU
$s4Pic2.Top    = 86.39992141784309     # 1097279 EMU
$s4Pic2.Width  = 305.26551818895615    # 3876872 EMU
$s4Pic2.Height = 432.0                 # 5486400 EMU

Write-Output "done"
